$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Test-result cells (C:G for rows 2-4) are stored as text (e.g. "544.0"),
# not numbers, in the source workbook. A leading apostrophe forces Excel
# to keep each replacement value as literal text instead of inferring a
# numeric type, matching the original "Test 1".."Test 5" column data.
$ws.Range("C2").Value = "'556.2"
$ws.Range("D2").Value = "'559.6"
$ws.Range("E2").Value = "'557.1"
$ws.Range("F2").Value = "'562.0"
$ws.Range("G2").Value = "'541.3"

$ws.Range("C3").Value = "'363.6"
$ws.Range("D3").Value = "'345.1"
$ws.Range("E3").Value = "'356.4"
$ws.Range("F3").Value = "'353.3"
$ws.Range("G3").Value = "'362.8"

$ws.Range("C4").Value = "'204.5"
$ws.Range("D4").Value = "'202.7"
$ws.Range("E4").Value = "'202.4"
$ws.Range("F4").Value = "'199.5"
$ws.Range("G4").Value = "'194.6"
